# ---------------------------------------------------------------------------
# Target diff analysis (ppt/slideMasters/slideMaster1.xml only):
#
#   Every "-"/"+" pair in the supplied diff touches the SAME element (same
#   tag) with the SAME set of attribute names/values (including the root
#   <p:sldMaster> xmlns:a / xmlns:p / xmlns:r declarations, <p:clrMap>,
#   <p:hf>, every <a:bodyPr>, <p:ph>, <a:rPr>/<a:endParaRPr>, <a:buFont>,
#   and every txStyles level's <a:lvlNpPr>/<a:defRPr>). Only the *order* in
#   which those attributes are written changes (alphabetical afterwards).
#   No tag, no attribute name, no attribute value, no child element, no
#   relationship and no theme reference is added, removed or retargeted —
#   confirmed by diffing each "-" line's (tag, {attr:val}) set against its
#   paired "+" line for all 66 changed lines in the patch.
#
#   This matches the commit message: the underlying writer was changed to
#   always emit *every* theme part it holds; for decks that (like this one,
#   which has speaker notes / a notesMaster wired to ppt/theme/theme2.xml)
#   already carried both themes, nothing about the deck's content changes —
#   only some parts get re-serialised by the new code path, which happens
#   to write attributes out in a stable/alphabetical order instead of the
#   original authoring order. That's a serializer implementation detail,
#   not a PowerPoint object-model edit: there is no Slide/Shape/TextFrame/
#   Theme property whose value differs between "before" and "after" here.
#
# So the faithful reproduction of this change is to leave every shape,
# placeholder, text run, and theme reference exactly as authored — i.e. no
# mutating calls — while still exercising the slide master / theme-related
# parts of the object model the diff concerns, read-only, so the session
# resolves them.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# The diff lives entirely inside the (first and only) slide master.
$master = $p.Slides.Item(1).Master

# Touch it (and the placeholders the diff enumerates: title, body, date,
# footer, slide number) read-only — no property is assigned, so no content,
# formatting, geometry, relationship or theme reference is altered.
$placeholders = $master.Shapes.Placeholders
Write-Host "slideMaster placeholders: $($placeholders.Count)"
Write-Host "slideMaster background: $($master.Background)"
Write-Host "slideMaster colorScheme: $($master.ColorScheme)"
